# The two Pearson/Edexcel logo pictures (in footer1.xml / footer2.xml)
# are renamed from "image2.png" to "image1.png", and the BTec logo
# picture (in header2.xml) is renamed from "image1.jpg" to "image2.jpg".
# Word's InlineShape object model has no settable "Name" property for
# inline pictures (that only exists on floating Shape/ShapeRange), so
# the rename of the wp:docPr/@name and pic:cNvPr/@name attributes is
# done by round-tripping the package through Document.WordOpenXML and
# doing a targeted text substitution there.

$d = $word.ActiveDocument

$xml = $d.WordOpenXML

$xml = $xml.Replace('name="image2.png"', 'name="image1.png"')
$xml = $xml.Replace('name="image1.jpg"', 'name="image2.jpg"')

$d.WordOpenXML = $xml
